$d = $word.ActiveDocument

function Find-ParaIndex($pattern) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs.Item($i).Range.Text -like $pattern) {
            return $i
        }
    }
    return -1
}

# 1. Delete the bold 36pt paragraph "GitHub is distributed version control system
#    that allows developer/programmer store their code in the cloud and collaborate."
#    entirely.
$introIdx = Find-ParaIndex("GitHub is distributed version control system*")
$d.Paragraphs.Item($introIdx).Range.Delete()

# 2. Change "Install GitHub : git-scm.com/download" to "git status".
$d.Content.Find.Execute("Install GitHub : git-scm.com/download", $false, $false, $false, $false, $false, $true, 1, $false, "git status", 2)

# 3. Remove every paragraph from "git init project1 (create local repository)"
#    through "Test GitHub" (including the trailing blank paragraphs), leaving only
#    the paragraph that used to read "launch command prompt ( git -version)"
#    (which becomes "git add .") followed directly by the final bookmark paragraph.
$startIdx = Find-ParaIndex("git init project1*")
$endIdx = Find-ParaIndex("Test GitHub*")
$startP = $d.Paragraphs.Item($startIdx)
$endP = $d.Paragraphs.Item($endIdx)
$d.Range($startP.Range.Start, $endP.Range.End).Delete()

# 4. Collapse the "launch command prompt ( git -version)" paragraph down to a
#    single run reading "git add ." (keeps the first run's formatting).
$d.Content.Find.Execute("launch command prompt ( git –version)", $false, $false, $false, $false, $false, $true, 1, $false, "git add .", 2)
